$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.974.11"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.395.38"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.81"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.65"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.394.02"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.974.03"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.373.68"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.46"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.039.48"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.21"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.50"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.19"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000130"
$ws.Range("E24").Value = "  +8.67%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.539.24"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.63"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.67"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.19"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.57"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.425.61"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.41"
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.91"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.11"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0790"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.78"
$ws.Range("E42").Value = "  +9.05%  "
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.788"
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.21"
$ws.Range("E46").Value = "  +7.06%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.36"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.91"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.18"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.348.28"
$ws.Range("E51").Value = "  +6.71%  "
